# Shift the Unintended_Deviation dataset forward by 2 days (retraining refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
$lastRow = $ws.Cells.Item(1, 1).Worksheet.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $tsCell = $ws.Cells.Item($r, 1)
    $quarterCell = $ws.Cells.Item($r, 4)
    $lookupCell = $ws.Cells.Item($r, 5)

    $serial = $tsCell.Value2
    $newSerial = $serial + 2
    $tsCell.Value2 = $newSerial

    $quarter = $quarterCell.Value2

    $newDate = $epoch.AddDays($newSerial)
    $dateStr = $newDate.ToString("dd.MM.yyyy")

    $lookupCell.Value = "$dateStr$quarter"
}
